$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new formula cell in row 1 (converts miles in R1 to rounded km)
$ws.Range("S1").Formula = '=ROUND(CONVERT(R1,"mi","km"),0)'

# Update the view: scroll so column P is the left-most visible column,
# and move the active selection to AC6
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 16
$ws.Range("AC6").Select()
